$d = $word.ActiveDocument

$replacements = @(
    @("850×9=", "615×8="),
    @("417×4=", "991×2="),
    @("416×6=", "405×2="),
    @("865×5=", "202×6="),
    @("846×4=", "710×6="),
    @("113×2=", "454×2="),
    @("919×7=", "432×4="),
    @("717×3=", "439×8="),
    @("327×7=", "444×3="),
    @("719×9=", "125×9="),
    @("776×2=", "609×6="),
    @("828×5=", "371×6="),
    @("566×3=", "742×4="),
    @("147×2=", "809×5="),
    @("328×6=", "269×4="),
    @("418×4=", "614×9="),
    @("849×3=", "597×5="),
    @("463×2=", "552×4="),
    @("495×5=", "158×8="),
    @("592×3=", "218×2="),
    @("716×2=", "510×2="),
    @("704×9=", "455×6="),
    @("317×9=", "207×9="),
    @("788×7=", "172×5="),
    @("396×5=", "246×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
